$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new observation for 2026/01/27 (19 o'clock) was recorded; it belongs
# chronologically right after the existing 2026/01/27 rows and before the
# 2026/12/29 block, so insert a fresh row at 725 and push everything below
# it down by one (shifts 725-766 -> 726-767).
$ws.Rows("725:725").Insert()

# Copy the date text from the row above so the new cell is stored as the
# same literal text "2026/01/27" (matches existing column-A formatting)
# instead of being auto-parsed into a date serial number.
$ws.Range("A724").Copy() | Out-Null
$ws.Range("A725").PasteSpecial(-4163)

$ws.Range("B725").Value = "火"
$ws.Range("C725").Value = 19
$ws.Range("D725").Value = 201

Write-Output "inserted row 725 for 2026/01/27"
